$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.200034379959106
$ws.Range("B1").Value = 2.06397008895874
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.101459980010986
$ws.Range("E1").Value = 1.207412600517273
